# api-planning.xlsx: add "Link" and "Database" sheets (MERN Project 31 -
# wiring the MongoDB connection info + a couple of setup/tutorial links),
# and move the active tab/selection onto the new "Database" sheet.

$wb = $excel.ActiveWorkbook

$apiSheet  = $wb.Worksheets.Item(1)
$respSheet = $wb.Worksheets.Item(2)

# --- New sheet: "Link" (right after "Response Status ") ---------------
$linkSheet = $wb.Worksheets.Add($null, $respSheet)
$linkSheet.Name = "Link"
$linkSheet.Range("A1").Value = "Tutorial"
$linkSheet.Range("A2").Value = "https://www.youtube.com/watch?v=nGWP07CtEx0&list=PLU4DS8KR-LJ0-MT2QfV-fvJiNorsoFs74&index=30"
$linkSheet.Range("A3").Value = "https://www.youtube.com/watch?v=413C1PlYIko"
$linkSheet.Range("H3").Select()

# --- New sheet: "Database" (right after "Link") ------------------------
$dbSheet = $wb.Worksheets.Add($null, $linkSheet)
$dbSheet.Name = "Database"
$dbSheet.Range("A1").Value = "ceepei14"
$dbSheet.Range("B1").Value = "hXvpMZrePqHSp2Yl"
$dbSheet.Range("B1").Select()

# --- Update the old selection on "API" (was J5, now F8) and drop its ---
# --- tabSelected flag by making "Database" the active sheet ------------
$apiSheet.Range("F8").Select()

$dbSheet.Activate()
